# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (fund position snapshot) between the
# existing "总计" (summary) and "2022-Q1" sheets, and records the new
# quarter's totals on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: insert the new sheet right after "总计" (i.e. before "2022-Q1")
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# A never-touched, still-default-styled cell on the summary sheet. Used
# purely as a "format donor" later to strip the automatic quote-prefix
# style Excel applies when a digit-only string is stored as text.
$cleanCell = $totalSheet.Cells.Item(1, 26)

# ---------------------------------------------------------------------
# Step 2: update "总计" — push the existing 2022-Q1 row down to row 3
# and write the new 2022-Q4 row into row 2
# ---------------------------------------------------------------------
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(3, 1).PasteSpecial($xlPasteFormats)

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(3, 3).Value = 1
$totalSheet.Cells.Item(3, 4).Value = 0.18

$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.02

# ---------------------------------------------------------------------
# Step 3: populate the new "2022-Q4" sheet with the fund holdings table
# ---------------------------------------------------------------------

# Header row (B1:H1) + first data column (A2:A3) reuse the bordered
# "title" style already used on row 1 / column A of "总计".
$totalSheet.Cells.Item(1, 2).Copy()
$q4Sheet.Cells.Item(1, 2).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(1, 3).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(1, 4).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(1, 5).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(1, 6).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(1, 7).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(1, 8).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(1, 2).Value = "基金代码"
$q4Sheet.Cells.Item(1, 3).Value = "基金名称"
$q4Sheet.Cells.Item(1, 4).Value = "基金规模"
$q4Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q4Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1, 8).Value = "仓位排名"

$totalSheet.Cells.Item(2, 1).Copy()
$q4Sheet.Cells.Item(2, 1).PasteSpecial($xlPasteFormats)
$q4Sheet.Cells.Item(3, 1).PasteSpecial($xlPasteFormats)

# Row 2 — 002872 华夏智胜价值成长股票C
$q4Sheet.Cells.Item(2, 1).Value = 0

$q4Sheet.Cells.Item(2, 2).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 2).Value = "002872"
$cleanCell.Copy()
$q4Sheet.Cells.Item(2, 2).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(2, 3).Value = "华夏智胜价值成长股票C"

$q4Sheet.Cells.Item(2, 4).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 4).Value = "2.68"
$cleanCell.Copy()
$q4Sheet.Cells.Item(2, 4).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(2, 5).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 5).Value = "92.97"
$cleanCell.Copy()
$q4Sheet.Cells.Item(2, 5).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(2, 6).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 6).Value = "0.64"
$cleanCell.Copy()
$q4Sheet.Cells.Item(2, 6).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(2, 7).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 7).Value = "0.0172"
$cleanCell.Copy()
$q4Sheet.Cells.Item(2, 7).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(2, 8).Value = 10

# Row 3 — 002871 华夏智胜价值成长股票A
$q4Sheet.Cells.Item(3, 1).Value = 1

$q4Sheet.Cells.Item(3, 2).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 2).Value = "002871"
$cleanCell.Copy()
$q4Sheet.Cells.Item(3, 2).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(3, 3).Value = "华夏智胜价值成长股票A"

$q4Sheet.Cells.Item(3, 4).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 4).Value = "0.92"
$cleanCell.Copy()
$q4Sheet.Cells.Item(3, 4).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(3, 5).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 5).Value = "92.97"
$cleanCell.Copy()
$q4Sheet.Cells.Item(3, 5).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(3, 6).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 6).Value = "0.64"
$cleanCell.Copy()
$q4Sheet.Cells.Item(3, 6).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(3, 7).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 7).Value = "0.0059"
$cleanCell.Copy()
$q4Sheet.Cells.Item(3, 7).PasteSpecial($xlPasteFormats)

$q4Sheet.Cells.Item(3, 8).Value = 10
